# Applies the LPD.docx edit described by the commit:
#  1. Remove the trailing bullet line about "Including a file in a patch..."
#     from the "4. WORKFLOW RULES / PRACTICES" paragraph.
#  2. Replace the "Decisions Added (2026-01-26)" heading + its four
#     Decision/Rationale/Implications/Revisit-Conditions paragraphs with a
#     new "Update - 2026-01-26" block containing four bolded
#     decision/note headers each followed by an indented explanation
#     paragraph.

$d = $word.ActiveDocument

# Helpful characters that don't type cleanly in plain ASCII source.
$bullet    = [char]0x2022
$enDash    = [char]0x2013
$lCurlyDbl = [char]0x201C
$rCurlyDbl = [char]0x201D

# ---------------------------------------------------------------------
# Change 1: delete the last bullet line in the WORKFLOW RULES paragraph.
# ---------------------------------------------------------------------
$oldBullet = $bullet + " Including a file in a patch implies intentional changes; unchanged files must not be included (even if re-saved, line-ending-normalized, or auto-formatted)."
$rng = $d.Content
$found = $rng.Find.Execute($oldBullet, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Delete()
}

# ---------------------------------------------------------------------
# Change 2: replace the "Decisions Added (2026-01-26)" section.
# ---------------------------------------------------------------------

# Locate the heading paragraph and the three paragraphs that follow it
# (there are 5 paragraphs total in the old section: heading, Decision,
# Rationale, Implications, Revisit Conditions).
$headingText = "Decisions Added (2026-01-26)"
$sectionStart = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $headingText) {
        $sectionStart = $i
        break
    }
}

if ($sectionStart -ne -1) {
    $pFirst = $d.Paragraphs.Item($sectionStart)
    $pLast  = $d.Paragraphs.Item($sectionStart + 4)
    $killRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)
    $killRange.Delete()
}

# Anchor = the paragraph now immediately before the insertion point
# (the "Revisit Conditions ... MAX_RECENTS_DISPLAY" paragraph).
$anchorIndex = $sectionStart - 1

# Index-based paragraph builder: after each InsertParagraphAfter call the
# newly created paragraph is the one whose index equals
# (previous paragraph index + 1); track that index manually.
$curIndex = $anchorIndex

function New-Para([string]$text, [bool]$bold, [bool]$indent, [int]$fontSize) {
    $curIndex = $curIndex + 1
    $prevIndex = $curIndex - 1
    $prevPara = $d.Paragraphs.Item($prevIndex)
    $prevPara.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($curIndex)
    $newRange = $newPara.Range
    if ($text -ne "") {
        $newRange.InsertBefore($text)
        # Refresh handles after text insertion (range bounds changed).
        $newPara = $d.Paragraphs.Item($curIndex)
        $newRange = $newPara.Range
    }
    if ($indent) {
        $newRange.ParagraphFormat.LeftIndent = 18
    }
    if ($bold -or $fontSize -gt 0) {
        $textRange = $d.Range($newRange.Start, $newRange.End - 1)
        if ($bold) {
            $textRange.Bold = 1
        }
        if ($fontSize -gt 0) {
            $textRange.Font.Size = $fontSize
        }
    }
}

# 1) blank separator paragraph
New-Para "" $false $false 0

# 2) "Update - 2026-01-26" (bold, 14pt)
$updateHeading = "Update " + $enDash + " 2026-01-26"
New-Para $updateHeading $true $false 14

# 3) Decision: Recents listbox scrolling fix (flex shrink).
New-Para "Decision: Recents listbox scrolling fix (flex shrink)." $true $false 0

# 4) body paragraph (indented)
$body1 = "Recents was not showing an internal scrollbar because listboxes were locked to computed heights (minHeight=height=maxHeight). In Office webviews the available height can be slightly smaller than computed, causing clipping instead of internal overflow. Fix: allow Favorites/Recents listboxes to shrink (minHeight:0, flex-shrink) while keeping overflowY:auto."
New-Para $body1 $false $true 0

# 5) Decision: Remove redundant "Prioritize Favorites" option.
$decision2 = "Decision: Remove redundant " + $lCurlyDbl + "Prioritize Favorites" + $rCurlyDbl + " option."
New-Para $decision2 $true $false 0

# 6) body paragraph (indented)
$body2 = "After finalizing the allocation rules (20% min / 80% cap + surplus donation), " + $lCurlyDbl + "Prioritize Favorites" + $rCurlyDbl + " is logically equivalent to setting the slider to the Favorites-max position (80/20). Kept the slider only."
New-Para $body2 $false $true 0

# 7) UI clarity: Slider label wording.
New-Para "UI clarity: Slider label wording." $true $false 0

# 8) body paragraph (indented)
$body3 = "Reword slider label to: " + $lCurlyDbl + "When space is limited, give more room to:" + $rCurlyDbl + " to communicate conditional effect in plain language."
New-Para $body3 $false $true 0

# 9) Workflow: Patch discipline and documentation.
New-Para "Workflow: Patch discipline and documentation." $true $false 0

# 10) body paragraph (indented)
$body4 = "Patches must (a) include an updated LPD whenever decisions/workflow rules change, (b) include a git commit message with " + $lCurlyDbl + "what was done before" + $rCurlyDbl + " context, and (c) include only files that actually changed (no unchanged files for context)."
New-Para $body4 $false $true 0
